# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (B15:J22) is re-sorted from descending period
# order (1910 -> 1904) to ascending period order (1904 -> 1910). The
# "Valor Mora" (column F) travels together with its period, so the
# period/value pairs themselves are unchanged - only their row order is.
#
# Before (row -> period : value):
#   16 -> 1910 : 42000
#   17 -> 1909 : 60000
#   18 -> 1908 : 60000
#   19 -> 1907 : 60000
#   20 -> 1906 : 33125
#   21 -> 1905 : 33125
#   22 -> 1904 : 33125
#
# After (row -> period : value):
#   16 -> 1904 : 33125
#   17 -> 1905 : 33125
#   18 -> 1906 : 33125
#   19 -> 1907 : 60000
#   20 -> 1908 : 60000
#   21 -> 1909 : 60000
#   22 -> 1910 : 42000

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("1904", "1905", "1906", "1907", "1908", "1909", "1910")
$valores  = @(33125, 33125, 33125, 60000, 60000, 60000, 42000)

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
